$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.22802204078956834
$ws.Range("A2").Value = -0.0059999999814408511
$ws.Range("A3").Value = -0.003999999983719249
$ws.Range("A4").Value = -0.007999999970534688
$ws.Range("A5").Value = -0.0029999999838885572
$ws.Range("A6").Value = -0.0019999999834876547
$ws.Range("A7").Value = -0.0099999999595117295
$ws.Range("A8").Value = -0.0099999999583291199
$ws.Range("A9").Value = -0.0019999999806139535
$ws.Range("A10").Value = -0.0019999999794837464
$ws.Range("A11").Value = -0.0029999999768399732
$ws.Range("A12").Value = -0.020322428029923234
$ws.Range("A13").Value = -0.003499999975599799
$ws.Range("A14").Value = -0.0079999999626396701
$ws.Range("A15").Value = -0.00099999998375999155
$ws.Range("A16").Value = 0.025306705869925672
$ws.Range("A17").Value = -0.0019999999810922375
$ws.Range("A18").Value = -0.003999999975077273
$ws.Range("A19").Value = -0.0039999999866799918
$ws.Range("A20").Value = -0.066737003933084793
$ws.Range("A21").Value = -0.0039999999801878516
$ws.Range("A22").Value = -0.0039999999800564012
$ws.Range("A23").Value = -0.0049999999794279049
$ws.Range("A24").Value = -0.032674708525368246
$ws.Range("A25").Value = -0.019999999931815005
$ws.Range("A26").Value = -0.0024999999818593466
$ws.Range("A27").Value = -0.0024999999815866758
$ws.Range("A28").Value = -0.0019999999819866332
$ws.Range("A29").Value = -0.006999999966586401
$ws.Range("A30").Value = 0.012106449437590427
$ws.Range("A31").Value = -0.0069999999652399225
$ws.Range("A32").Value = -0.0099999999561894981
$ws.Range("A33").Value = -0.003999999973768098
